# Fruta / hortaliza, semanal
# Insert two new weekly price rows (new row 441 and 442) at the top of the
# "Alcachofa" data block on the active sheet. Every existing data row from
# 441..502 shifts down by two (to 443..504); their content is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 441 and 442 downward, inserting two blank rows in their place.
$ws.Range("441:442").Insert()

# --- New row 441 ---------------------------------------------------------
$ws.Range("A441").Value = 3
$ws.Range("B441").Value = "Femacal de La Calera"
$ws.Range("C441").Value = "Coquimbo"
$ws.Range("D441").Value = 45077
$ws.Range("E441").Value = 5
$ws.Range("F441").Value = 100112013
$ws.Range("G441").Value = "Alcachofa"
$ws.Range("H441").Value = "Argentina(o)"
$ws.Range("I441").Value = "Primera"
$ws.Range("J441").Value = 90
$ws.Range("K441").Value = 13000
$ws.Range("L441").Value = 13500
$ws.Range("M441").Value = 13250
$ws.Range("N441").Value = "`$/caja 50 unidades"
$ws.Range("O441").Value = "Provincia de Limarí"
$ws.Range("P441").Value = 265
$ws.Range("Q441").Value = 50
$ws.Range("R441").Value = "Hortaliza"

# --- New row 442 ---------------------------------------------------------
$ws.Range("A442").Value = 3
$ws.Range("B442").Value = "Femacal de La Calera"
$ws.Range("C442").Value = "Coquimbo"
$ws.Range("D442").Value = 45077
$ws.Range("E442").Value = 5
$ws.Range("F442").Value = 100112013
$ws.Range("G442").Value = "Alcachofa"
$ws.Range("H442").Value = "Española"
$ws.Range("I442").Value = "Primera"
$ws.Range("J442").Value = 80
$ws.Range("K442").Value = 16000
$ws.Range("L442").Value = 16000
$ws.Range("M442").Value = 16000
$ws.Range("N442").Value = "`$/caja 30 unidades"
$ws.Range("O442").Value = "Provincia de Limarí"
$ws.Range("P442").Value = 533
$ws.Range("Q442").Value = 30
$ws.Range("R442").Value = "Hortaliza"
